# Apply the per-cell "Price" (D) / "Volume(1h)" (E) updates, plus the
# three coin-row re-sorts (B/C/D/E on rows 36-37 and 40-43), exactly as
# captured by the source OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.083.71"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.652.36"
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'217.40"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'0.5246"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.2600"
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("D9").Value = "'0.06325"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").Value = "'20.37"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("D11").Value = "'0.07790"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "'4.509"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "1.651.00"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").Value = "1.879.01"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "'0.5500"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "0.0₅8212"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "'65.56"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "26.097.10"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'4.587"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Value = "'190.79"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").Value = "'10.08"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'6.040"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'144.04"
$ws.Range("E25").Value = "  +3.24%  "
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").Value = "'7.239"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'16.06"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").Value = "'0.05836"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "'3.547"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'3.266"
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").Value = "'1.586"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").Value = "'0.9464"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "'2.780"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.410"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "'0.5742"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").Value = "'0.01608"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.748"
$ws.Range("E40").Value = "  -4.85%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'104.13"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8417"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("D44").Value = "1.030.62"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").Value = "1.794.72"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "'57.12"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").Value = "'0.4328"
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D50").Value = "'7.816"
$ws.Range("E50").Value = "  -3.07%  "
$ws.Range("D51").Value = "'1.457"
$ws.Range("E51").Value = "  +0.32%  "
